# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Map event name (column C) -> new value for column F
    $updates = @{
        "南昌·Sunflower Garden动漫游戏展" = 6340
        "赣州·卡尼动漫展" = 369
        "鹰潭·MZD动漫游戏嘉年华" = 53
        "赣州·十万伏特-第七届青年文化综合展览会" = 56
        "景德镇·第十六届瓷都ACG动漫游戏博览会" = 614
        "南昌·萌卡动漫展" = 3087
        "景德镇·第十六届瓷都ACG内场—花玲&宴宁" = 180
        "江西·JMG（江西广电）第二届UP动漫游戏博览会" = 1734
    }

    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        if ($name -ne $null -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
